$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mapping of row -> (DAMSLTag, DialogAct) updates for columns I and J
$updates = @(
    @{Row=2; I="sd"; J="Statement-non-opinion"}
    @{Row=7; I="sd"; J="Statement-non-opinion"}
    @{Row=14; I="sv"; J="Statement-opinion"}
    @{Row=28; I="sd"; J="Statement-non-opinion"}
    @{Row=54; I="sv"; J="Statement-opinion"}
    @{Row=55; I="sv"; J="Statement-opinion"}
    @{Row=57; I="sv"; J="Statement-opinion"}
    @{Row=65; I="sd"; J="Statement-non-opinion"}
    @{Row=72; I="sd"; J="Statement-non-opinion"}
    @{Row=84; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=85; I="sd"; J="Statement-non-opinion"}
    @{Row=86; I="sd"; J="Statement-non-opinion"}
    @{Row=90; I="sd"; J="Statement-non-opinion"}
    @{Row=91; I="sd"; J="Statement-non-opinion"}
    @{Row=94; I="ba"; J="Appreciation"}
    @{Row=112; I="sd"; J="Statement-non-opinion"}
    @{Row=113; I="sd"; J="Statement-non-opinion"}
    @{Row=121; I="aa"; J="Agree/Accept"}
    @{Row=125; I="aa"; J="Agree/Accept"}
    @{Row=128; I="aa"; J="Agree/Accept"}
    @{Row=130; I="%"; J="Uninterpretable"}
    @{Row=143; I="sd"; J="Statement-non-opinion"}
    @{Row=155; I="sd"; J="Statement-non-opinion"}
    @{Row=156; I="sd"; J="Statement-non-opinion"}
    @{Row=159; I="sd"; J="Statement-non-opinion"}
    @{Row=164; I="sd"; J="Statement-non-opinion"}
    @{Row=175; I="%"; J="Uninterpretable"}
    @{Row=176; I="aa"; J="Agree/Accept"}
    @{Row=203; I="sd"; J="Statement-non-opinion"}
    @{Row=206; I="sv"; J="Statement-opinion"}
    @{Row=211; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=231; I="aa"; J="Agree/Accept"}
    @{Row=235; I="%"; J="Uninterpretable"}
    @{Row=238; I="sd"; J="Statement-non-opinion"}
    @{Row=257; I="sd"; J="Statement-non-opinion"}
    @{Row=261; I="sd"; J="Statement-non-opinion"}
    @{Row=266; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=267; I="sd"; J="Statement-non-opinion"}
    @{Row=280; I="sv"; J="Statement-opinion"}
    @{Row=282; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=286; I="sv"; J="Statement-opinion"}
    @{Row=295; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=302; I="sd"; J="Statement-non-opinion"}
    @{Row=306; I="sv"; J="Statement-opinion"}
    @{Row=311; I="sd"; J="Statement-non-opinion"}
    @{Row=317; I="aa"; J="Agree/Accept"}
    @{Row=319; I="sd"; J="Statement-non-opinion"}
    @{Row=324; I="qy"; J="Yes-No-Question"}
    @{Row=332; I="sd"; J="Statement-non-opinion"}
    @{Row=334; I="aa"; J="Agree/Accept"}
    @{Row=350; I="sd"; J="Statement-non-opinion"}
    @{Row=356; I="qy"; J="Yes-No-Question"}
    @{Row=357; I="sd"; J="Statement-non-opinion"}
    @{Row=359; I="%"; J="Uninterpretable"}
    @{Row=360; I="aa"; J="Agree/Accept"}
    @{Row=363; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=364; I="%"; J="Uninterpretable"}
    @{Row=377; I="%"; J="Uninterpretable"}
    @{Row=383; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=385; I="aa"; J="Agree/Accept"}
    @{Row=394; I="aa"; J="Agree/Accept"}
    @{Row=403; I="aa"; J="Agree/Accept"}
    @{Row=404; I="aa"; J="Agree/Accept"}
    @{Row=411; I="ba"; J="Appreciation"}
    @{Row=412; I="sd"; J="Statement-non-opinion"}
    @{Row=425; I="sd"; J="Statement-non-opinion"}
    @{Row=427; I="sd"; J="Statement-non-opinion"}
    @{Row=454; I="aa"; J="Agree/Accept"}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}

$wb.Save()
